# "Enterprises density (per 1000 people)" (value 55) moves up to sit right
# after the "MSMEs" header, ahead of "Employment (% of total)" and
# "Enterprises (absolute #)" - those two rows shift down to make room:
#
#   before                                   after
#   12: Employment (% of total) | 86.2       12: Enterprises density (per 1000 people) | 55
#   13: Enterprises (absolute #) | 1279784   13: Employment (% of total) | 86.2
#   14: Enterprises density (...) | 55       14: Enterprises (absolute #) | 1279784

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(12, 1).Value = "Enterprises density (per 1000 people)"
$ws.Cells.Item(13, 1).Value = "Employment (% of total)"
$ws.Cells.Item(14, 1).Value = "Enterprises (absolute #)"

# The D column values look numeric ("55", "86.2", "1279784") but must stay
# stored as text, matching the original cells. Assigning them straight to
# .Value would auto-coerce to numbers, so build each as a text formula in
# an unused scratch cell and paste its value back in - this keeps the
# destination cell's original (text) type and style untouched.
$scratch = $ws.Cells.Item(1, 10)

$scratch.Formula = "=""55"""
$scratch.Copy()
$ws.Cells.Item(12, 4).PasteSpecial(-4163)

$scratch.Formula = "=""86.2"""
$scratch.Copy()
$ws.Cells.Item(13, 4).PasteSpecial(-4163)

$scratch.Formula = "=""1279784"""
$scratch.Copy()
$ws.Cells.Item(14, 4).PasteSpecial(-4163)

$scratch.Clear()
